$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> new cell text, per the diff. Word rows are 1-indexed;
# only rows 1, 5, 10, 15, 20 hold data (the remaining rows are blank spacer rows).
$updates = @(
    @{ Row = 1;  Col = 1; Old = "480×3=1440"; New = "520×2=1040" },
    @{ Row = 1;  Col = 2; Old = "168×6=1008"; New = "628×6=3768" },
    @{ Row = 1;  Col = 3; Old = "699×4=2796"; New = "155×8=1240" },
    @{ Row = 1;  Col = 4; Old = "209×2=418";  New = "980×3=2940" },
    @{ Row = 1;  Col = 5; Old = "996×7=6972"; New = "712×2=1424" },

    @{ Row = 5;  Col = 1; Old = "430×5=2150"; New = "287×2=574" },
    @{ Row = 5;  Col = 2; Old = "149×8=1192"; New = "963×7=6741" },
    @{ Row = 5;  Col = 3; Old = "895×3=2685"; New = "529×7=3703" },
    @{ Row = 5;  Col = 4; Old = "366×9=3294"; New = "626×6=3756" },
    @{ Row = 5;  Col = 5; Old = "517×2=1034"; New = "688×7=4816" },

    @{ Row = 10; Col = 1; Old = "976×9=8784"; New = "727×9=6543" },
    @{ Row = 10; Col = 2; Old = "671×7=4697"; New = "772×5=3860" },
    @{ Row = 10; Col = 3; Old = "371×7=2597"; New = "431×6=2586" },
    @{ Row = 10; Col = 4; Old = "559×6=3354"; New = "282×4=1128" },
    @{ Row = 10; Col = 5; Old = "241×3=723";  New = "425×6=2550" },

    @{ Row = 15; Col = 1; Old = "761×4=3044"; New = "805×7=5635" },
    @{ Row = 15; Col = 2; Old = "321×4=1284"; New = "818×7=5726" },
    @{ Row = 15; Col = 3; Old = "574×9=5166"; New = "249×2=498" },
    @{ Row = 15; Col = 4; Old = "723×6=4338"; New = "869×8=6952" },
    @{ Row = 15; Col = 5; Old = "292×4=1168"; New = "486×5=2430" },

    @{ Row = 20; Col = 1; Old = "451×8=3608"; New = "853×6=5118" },
    @{ Row = 20; Col = 2; Old = "464×3=1392"; New = "878×4=3512" },
    @{ Row = 20; Col = 3; Old = "256×6=1536"; New = "532×3=1596" },
    @{ Row = 20; Col = 4; Old = "644×5=3220"; New = "310×9=2790" },
    @{ Row = 20; Col = 5; Old = "517×2=1034"; New = "862×8=6896" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Shrink the range so Find doesn't touch the cell-end/paragraph marks.
    $rng.End = $rng.End - 1
    $found = $rng.Find.Execute($u.Old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $u.New, 2)
    if (-not $found) {
        Write-Host "WARNING: not found Row=$($u.Row) Col=$($u.Col) Old=$($u.Old)"
    }
}

Write-Host "Done"
